# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# right before column N ("Late"), pushing the existing N/O/P columns
# (Late / heading / Outstanding) one position to the right, and the
# sheet becomes the active/selected tab with the cursor left on R10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet/tab (was "Transactions").
$ws.Activate()

# Insert a brand-new, empty column at N - everything from N onward
# (N, O, P) shifts one column to the right (N->O, O->P, P->Q).
$ws.Columns("N").Insert()

# The newly inserted column picks up the same width as the column
# immediately to its left (M), same as Excel does by default.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Leave the selection where the user ended up after the edit.
$ws.Range("R10").Select() | Out-Null
